# Weekly data refresh: a new daily record is prepended to the Pomelo
# price block (Feria Lagunitas de Puerto Montt). Insert a new row at
# row 60 - this pushes the existing rows 60..144 down to 61..145 - and
# populate the new row with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60; everything below shifts down by one.
$ws.Rows(60).Insert()

# Fill in the new record at row 60.
$ws.Cells.Item(60, 1).Value  = 4
$ws.Cells.Item(60, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(60, 3).Value  = "Los Lagos"
$ws.Cells.Item(60, 4).Value2 = 44467
$ws.Cells.Item(60, 5).Value  = 10
$ws.Cells.Item(60, 6).Value  = "Fruta"
$ws.Cells.Item(60, 7).Value  = 100102
$ws.Cells.Item(60, 8).Value  = "Cítricos"
$ws.Cells.Item(60, 9).Value  = 100102006
$ws.Cells.Item(60, 10).Value = "Pomelo"
$ws.Cells.Item(60, 11).Value = "Start Ruby"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 160
$ws.Cells.Item(60, 14).Value = 12000
$ws.Cells.Item(60, 15).Value = 12000
$ws.Cells.Item(60, 16).Value = 12000
$ws.Cells.Item(60, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(60, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(60, 19).Value = 857
$ws.Cells.Item(60, 20).Value = 14
